$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Find the next empty row after the current data (row 58 -> new row 59)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# Copy the formatting (borders, fill, alignment, wrap text) from the row above
$srcRange = $ws.Range($ws.Cells.Item($lastRow, 1), $ws.Cells.Item($lastRow, 2))
$dstRange = $ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, 2))
$srcRange.Copy()
$dstRange.PasteSpecial(-4122)

# Set the new date label in column A
$ws.Cells.Item($newRow, 1).Value = "14-11-2025"

# Set the new gold price description in column B
$ws.Cells.Item($newRow, 2).Value = "The price of gold in India today is ₹12,785 per gram for 24 karat gold, ₹11,720 per gram for 22 karat gold and ₹9,589 per gram for 18 karat gold (also called 999 gold)."

$wb.Save()
